$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of A9 (border-only style) onto the new cell A10
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Add the new row content
$ws.Range("A10").Value = "SSH Key Added"

# Update the active selection to the newly added cell
$ws.Range("A10").Select()
